# Weekly update: insert a new record at row 74 (pushing existing rows down
# by one) for "Feria Lagunitas de Puerto Montt - Haba". The new row carries
# the same fixed attributes as the row it displaces, but with a new date
# and updated price/volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 74 (and everything below it) down by one row.
$ws.Range("A74").EntireRow.Insert()

# Populate the freshly inserted row 74 with the new weekly record.
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44754
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 100112026
$ws.Cells.Item(74, 7).Value = "Haba"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 90
$ws.Cells.Item(74, 11).Value = 25000
$ws.Cells.Item(74, 12).Value = 25000
$ws.Cells.Item(74, 13).Value = 25000
$ws.Cells.Item(74, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(74, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(74, 16).Value = 1000
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = "Hortaliza"
